$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-12 for columns I (I0) and J (IF)
$data = @{
    2  = @(1, 6)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 4)
    6  = @(1, 4)
    7  = @(9, 9)
    8  = @(1, 3)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(4, 7)
    12 = @(7, 8)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
